{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n//\n// Change summary (content-level):\n//   1. Turn the plain-text URL at the end of the \"Jim (2016).\" paragraph into\n//      a real hyperlink (same text/URL).\n//   2. Append three new \"label: <link>\" paragraphs citing additional\n//      resources (SFX/BG music tutorial, Royalty-free music, Royalty-fee\n//      sounds) right after the \"Jim (2016).\" paragraph.\n//   3. Append one trailing empty paragraph at the very end of the body.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\n// --- 1. Locate the \"Jim (2016).\" paragraph and hyperlink-ify its URL -------\nconst jimUrl =\n  \"https://stackoverflow.com/questions/34695396/how-to-render-a-paragraph-onto-a-surface-in-pygame\";\n\nlet jimParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(jimUrl) !== -1) {\n    jimParagraph = paragraphs.items[i];\n    break;\n  }\n}\nif (!jimParagraph) {\n  throw new Error('Could not find the \"Jim (2016).\" paragraph.');\n}\n\nconst jimUrlHits = jimParagraph.search(jimUrl, { matchCase: true });\njimUrlHits.load(\"items\");\nawait context.sync();\n\njimUrlHits.items[0].hyperlink = jimUrl;\nawait context.sync();\n\n// --- 2. Append the three new citation paragraphs ---------------------------\nfunction appendCitation(afterParagraph, label, url) {\n  const newParagraph = afterParagraph.insertParagraph(label + url, \"After\");\n  const urlHits = newParagraph.search(url, { matchCase: true });\n  urlHits.load(\"items\");\n  return { newParagraph, urlHits };\n}\n\nconst citations = [\n  {\n    label: \"SFX/BG music tutorial: \",\n    url: \"https://opensource.com/article/20/9/add-sound-python-game\",\n  },\n  {\n    label: \"Royalty-free music: \",\n    url: \"https://filmmusic.io/artists/kevin-macleod\",\n  },\n  {\n    label: \"Royalty-fee sounds: \",\n    url: \"https://freesound.org/\",\n  },\n];\n\nlet cursor = jimParagraph;\nconst pending = [];\nfor (const c of citations) {\n  const { newParagraph, urlHits } = appendCitation(cursor, c.label, c.url);\n  pending.push({ urlHits, url: c.url });\n  cursor = newParagraph;\n}\nawait context.sync();\n\nfor (const p of pending) {\n  p.urlHits.items[0].hyperlink = p.url;\n}\nawait context.sync();\n\n// --- 3. Trailing empty paragraph -------------------------------------------\ncursor.insertParagraph(\"\", \"After\");\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is the open document.\n#\n# Change summary (content-level):\n#   1. Turn the plain-text URL at the end of the \"Jim (2016).\" paragraph into\n#      a real hyperlink (same text/URL).\n#   2. Append three new \"label: <link>\" paragraphs citing additional\n#      resources (SFX/BG music tutorial, Royalty-free music, Royalty-fee\n#      sounds) right after the \"Jim (2016).\" paragraph.\n#   3. Append one trailing empty paragraph at the very end of the body.\n\n$d = $word.ActiveDocument\n\n# --- 1. Locate the \"Jim (2016).\" paragraph and hyperlink-ify its URL -------\n$jimUrl = \"https://stackoverflow.com/questions/34695396/how-to-render-a-paragraph-onto-a-surface-in-pygame\"\n\n$jimParagraph = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.Contains($jimUrl)) {\n        $jimParagraph = $p\n        break\n    }\n}\nif ($jimParagraph -eq $null) {\n    throw \"Could not find the 'Jim (2016).' paragraph.\"\n}\n\n$jimUrlRange = $jimParagraph.Range.Duplicate\n$jimUrlRange.Find.Execute($jimUrl) | Out-Null\n$d.Hyperlinks.Add($jimUrlRange, $jimUrl, \"\", \"\", $jimUrl) | Out-Null\n\n# --- 2. Append the three new citation paragraphs ---------------------------\nfunction Add-CitationParagraph($afterParagraph, $label, $url) {\n    $insertRange = $afterParagraph.Range\n    $insertRange.Collapse(0)   # wdCollapseEnd\n    $insertRange.InsertParagraphAfter()\n\n    $newParagraph = $d.Paragraphs.Item($afterParagraph.Index + 1)\n    $newRange = $newParagraph.Range\n    $newRange.Collapse(0)\n    $newRange.InsertAfter($label + $url)\n\n    $urlRange = $newParagraph.Range.Duplicate\n    $urlRange.Find.Execute($url) | Out-Null\n    $d.Hyperlinks.Add($urlRange, $url, \"\", \"\", $url) | Out-Null\n\n    return $newParagraph\n}\n\n$cursor = $jimParagraph\n$cursor = Add-CitationParagraph $cursor \"SFX/BG music tutorial: \" \"https://opensource.com/article/20/9/add-sound-python-game\"\n$cursor = Add-CitationParagraph $cursor \"Royalty-free music: \" \"https://filmmusic.io/artists/kevin-macleod\"\n$cursor = Add-CitationParagraph $cursor \"Royalty-fee sounds: \" \"https://freesound.org/\"\n\n# --- 3. Trailing empty paragraph -------------------------------------------\n$endRange = $cursor.Range\n$endRange.Collapse(0)\n$endRange.InsertParagraphAfter()\n"}
